$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3 content ---
$ws.Range("A3").Value = "change_request_2"
$ws.Range("B3").Value = "Team member change"
$ws.Range("C3").Value = "1-Future tasks assignation (NON)`n2--Running tasks assignation (NON)`n3-Training needed`n4-Training for new member on her old team (NON)"
$ws.Range("D3").Value = "1-Old team member and leader explained  plan, objective, documents, links and how we doing tasks for new member.`n-Expected time  (2 hrs)."

# --- Header row (row1): add vertical centering ---
$ws.Range("A1:D1").VerticalAlignment = -4108  # xlCenter

# --- Body style used by row2 (A:D) and row3 (A,C,D): bold 11, vertical center, wrap text ---
$ws.Range("A2:D2").Font.Bold = $true
$ws.Range("A2:D2").Font.Size = 11
$ws.Range("A2:D2").VerticalAlignment = -4108  # xlCenter
$ws.Range("A2:D2").WrapText = $true

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").VerticalAlignment = -4108  # xlCenter
$ws.Range("A3").WrapText = $true

$ws.Range("C3").Font.Bold = $true
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").VerticalAlignment = -4108  # xlCenter
$ws.Range("C3").WrapText = $true

$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Font.Size = 11
$ws.Range("D3").VerticalAlignment = -4108  # xlCenter
$ws.Range("D3").WrapText = $true

# --- Body style for B3: bold 11, vertical center, no wrap ---
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").VerticalAlignment = -4108  # xlCenter
$ws.Range("B3").WrapText = $false

# --- Row heights ---
$ws.Rows(2).RowHeight = 185.25
$ws.Rows(3).RowHeight = 72

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 24.140625
$ws.Columns("B").ColumnWidth = 39
$ws.Columns("C").ColumnWidth = 46.85546875
$ws.Columns("D").ColumnWidth = 39

# --- Selection ---
$ws.Range("E2").Select()
